$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.340.43"
$ws.Range("E2").Value = "  +0.14%  "
$ws.Range("D3").Value = "3.034.63"
$ws.Range("E3").Value = "  +1.05%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.27"
$ws.Range("E5").Value = "  -0.51%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.09"
$ws.Range("E6").Value = "  +2.94%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "3.031.43"
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.521"
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.65"
$ws.Range("E10").Value = "  +0.77%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("E11").Value = "  -1.46%  "
$ws.Range("E12").Value = "  +7.36%  "
$ws.Range("E13").Value = "  -1.44%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.69"
$ws.Range("E14").Value = "  +6.06%  "
$ws.Range("E15").Value = "  +0.04%  "
$ws.Range("D16").Value = "66.325.31"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("D17").Value = "3.536.34"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("E18").Value = "  +4.52%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "16.59"
$ws.Range("E19").Value = "  +19.68%  "
$ws.Range("D20").Value = "3.031.41"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "466.78"
$ws.Range("E21").Value = "  +3.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.711"
$ws.Range("E22").Value = "  +3.77%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.41"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.08"
$ws.Range("E24").Value = "  +0.88%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.75"
$ws.Range("E25").Value = "  +3.90%  "
$ws.Range("E26").Value = "  -1.38%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.05"
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.21"
$ws.Range("E29").Value = "  +1.03%  "
$ws.Range("E30").Value = "  +0.74%  "
$ws.Range("E31").Value = "  +1.37%  "
$ws.Range("E32").Value = "  +6.47%  "
$ws.Range("D33").Value = "0.0₃0993"
$ws.Range("E33").Value = "  -4.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "28.17"
$ws.Range("E34").Value = "  +3.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.992"
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("E37").Value = "  +0.63%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.42"
$ws.Range("E38").Value = "  +9.67%  "
$ws.Range("E39").Value = "  -0.92%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.316"
$ws.Range("E40").Value = "  +2.11%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "49.55"
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  +2.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.84"
$ws.Range("E44").Value = "  -3.91%  "
$ws.Range("E45").Value = "  +0.49%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "379.50"
$ws.Range("E46").Value = "  -5.78%  "
$ws.Range("D47").Value = "2.709.90"
$ws.Range("E47").Value = "  -2.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "134.15"
$ws.Range("E48").Value = "  +0.62%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.49"
$ws.Range("E50").Value = "  +2.75%  "
$ws.Range("E51").Value = "  +3.82%  "
